# Edit: PrecioFrutaHortalizas - Hortaliza, Vega Modelo de Temuco - Albahaca
# Commit message: "Fruta / hortaliza, semanal"
#
# The weekly update adds two new daily price records to the dataset.
# In the shipped workbook this shows up as:
#   - a new row inserted at row 303 (pushes the former rows 303..406 down by
#     one, to 304..407)
#   - a second new row inserted at (the now-shifted) row 407 (pushes the
#     former rows 406..407, now sitting at 407..408, further down to 408..409)
#   - <dimension> grows from A1:R407 to A1:R409
#
# Using native Rows(...).Insert() reproduces the shift (and carries the
# per-column number formatting, e.g. the date style on column D) without us
# having to rewrite the ~400 untouched rows by hand.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Insert the first new record at row 303.
# ---------------------------------------------------------------------------
$ws.Rows(303).Insert()

$ws.Cells.Item(303, 1).Value = 10
$ws.Cells.Item(303, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(303, 3).Value = "La Araucanía"
$ws.Cells.Item(303, 4).Value = 45120
$ws.Cells.Item(303, 5).Value = 9
$ws.Cells.Item(303, 6).Value = 100112052
$ws.Cells.Item(303, 7).Value = "Albahaca"
$ws.Cells.Item(303, 8).Value = "Sin especificar"
$ws.Cells.Item(303, 9).Value = "Primera"
$ws.Cells.Item(303, 10).Value = 150
$ws.Cells.Item(303, 11).Value = 5000
$ws.Cells.Item(303, 12).Value = 6000
$ws.Cells.Item(303, 13).Value = 5667
$ws.Cells.Item(303, 14).Value = "$/paquete"
$ws.Cells.Item(303, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(303, 16).Value = 5667
$ws.Cells.Item(303, 17).Value = 1
$ws.Cells.Item(303, 18).Value = "Hortaliza"

# ---------------------------------------------------------------------------
# 2) Insert the second new record. After step 1, the former row 406 sits at
#    row 407 and the former row 407 sits at row 408, so inserting here (at
#    row 407) places the new record between them, finally landing the two
#    old rows at 408 and 409.
# ---------------------------------------------------------------------------
$ws.Rows(407).Insert()

$ws.Cells.Item(407, 1).Value = 10
$ws.Cells.Item(407, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(407, 3).Value = "La Araucanía"
$ws.Cells.Item(407, 4).Value = 45121
$ws.Cells.Item(407, 5).Value = 9
$ws.Cells.Item(407, 6).Value = 100112052
$ws.Cells.Item(407, 7).Value = "Albahaca"
$ws.Cells.Item(407, 8).Value = "Sin especificar"
$ws.Cells.Item(407, 9).Value = "Primera"
$ws.Cells.Item(407, 10).Value = 80
$ws.Cells.Item(407, 11).Value = 6000
$ws.Cells.Item(407, 12).Value = 6000
$ws.Cells.Item(407, 13).Value = 6000
$ws.Cells.Item(407, 14).Value = "$/paquete"
$ws.Cells.Item(407, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(407, 16).Value = 6000
$ws.Cells.Item(407, 17).Value = 1
$ws.Cells.Item(407, 18).Value = "Hortaliza"
